# Update "want-to-go" counts (column F) on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 103
$ws1.Range("F3").Value = 4132
$ws1.Range("F6").Value = 17
$ws1.Range("F8").Value = 38
$ws1.Range("F9").Value = 204
$ws1.Range("F11").Value = 102
$ws1.Range("F12").Value = 145
$ws1.Range("F13").Value = 1547
$ws1.Range("F14").Value = 282
$ws1.Range("F15").Value = 3051
$ws1.Range("F16").Value = 210

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 103
$ws4.Range("F3").Value = 4132
$ws4.Range("F6").Value = 17
$ws4.Range("F9").Value = 38
$ws4.Range("F11").Value = 204
$ws4.Range("F13").Value = 102
$ws4.Range("F14").Value = 145
$ws4.Range("F17").Value = 1547
$ws4.Range("F18").Value = 282
$ws4.Range("F19").Value = 3051
$ws4.Range("F20").Value = 210
